# VerveStacks_CHE_grids model update - 2025-08-16 12:12
#
# 1. Rename the "varbl map" sheet to "timeslice map" and replace its
#    contents (which used to describe the ~Varbl_map dimension table,
#    now unused) with a new ~Timeslice_Map dimension table describing
#    ts_type / ts_season groupings.
# 2. Make "timeslice map" the active/selected sheet (it used to be
#    TS_Defs).

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheet and rebuild its contents -----------------------------
$ws = $wb.Worksheets.Item("varbl map")
$ws.Name = "timeslice map"

# Clear out the old ~Varbl_map table entirely before writing the new one.
$ws.Cells.Clear()

$ws.Range("A1").Value = "~Timeslice_Map"

$ws.Range("A2").Value = "dimension"
$ws.Range("B2").Value = "name"
$ws.Range("C2").Value = "description"

$ws.Range("A3").Value = "ts_type"
$ws.Range("B3").Value = "*,-s?a*"
$ws.Range("C3").Value = "hourly"

$ws.Range("A4").Value = "ts_type"
$ws.Range("B4").Value = "s?a*"
$ws.Range("C4").Value = "aggregated"

$ws.Range("A5").Value = "ts_season"
$ws.Range("B5").Value = "S1*"
$ws.Range("C5").Formula = "=LEFT(B5,2)"

$ws.Range("A6").Value = "ts_season"
$ws.Range("B6").Value = "S2*"
$ws.Range("C6").Formula = "=LEFT(B6,2)"

$ws.Range("A7").Value = "ts_season"
$ws.Range("B7").Value = "S3*"
$ws.Range("C7").Formula = "=LEFT(B7,2)"

$ws.Range("A8").Value = "ts_season"
$ws.Range("B8").Value = "S4*"
$ws.Range("C8").Formula = "=LEFT(B8,2)"

$ws.Range("A9").Value = "ts_season"
$ws.Range("B9").Value = "S5*"
$ws.Range("C9").Formula = "=LEFT(B9,2)"

$ws.Range("A10").Value = "ts_season"
$ws.Range("B10").Value = "S6*"
$ws.Range("C10").Formula = "=LEFT(B10,2)"

# --- 2. Make this sheet the active / selected tab -------------------------
$ws.Activate()
$ws.Range("A2").Select()
